$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with rich text runs) ---
# A8: "Volume 30   Number  31" -> "...32" (edit last run via Characters)
$a8 = $ws.Range("A8")
$a8text = $a8.Value2
$a8idx = $a8text.LastIndexOf("31") + 1
$a8.Characters($a8idx, 2).Text = "32"

# C9: "Report Covering the Week  7/31/2023  Through  8/6/2023"
#     -> "...8/7/2023  Through  8/13/2023"
$c9 = $ws.Range("C9")
$c9text = $c9.Value2
$c9idx1 = $c9text.IndexOf("7/31/2023") + 1
$c9.Characters($c9idx1, 9).Text = "8/7/2023"
$c9text2 = $c9.Value2
$c9idx2 = $c9text2.IndexOf("8/6/2023") + 1
$c9.Characters($c9idx2, 8).Text = "8/13/2023"

# --- Cells changing from numeric to shared-text "0" ---
# Reference cell already styled/typed as text "0": use D15 (style with numFmt General, text "0")
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Plain numeric value updates ---
$ws.Range("L15").Value = -50
$ws.Range("M15").Value = 14.285714285714
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -9.090909090909
$ws.Range("I16").Value = 101
$ws.Range("J16").Value = 117
$ws.Range("K16").Value = -13.675213675213
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = -6.481481481481
$ws.Range("N16").Value = -84.603658536585
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 233.333333333333
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 23.809523809523
$ws.Range("I17").Value = 140
$ws.Range("J17").Value = 149
$ws.Range("K17").Value = -6.040268456375
$ws.Range("L17").Value = 15.702479338843
$ws.Range("M17").Value = 77.215189873417
$ws.Range("N17").Value = -18.60465116279
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -45.833333333333
$ws.Range("I18").Value = 135
$ws.Range("J18").Value = 182
$ws.Range("K18").Value = -25.824175824175
$ws.Range("L18").Value = -2.877697841726
$ws.Range("M18").Value = -21.052631578947
$ws.Range("N18").Value = -86.620416253716
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -19.047619047619
$ws.Range("F19").Value = 85
$ws.Range("G19").Value = 80
$ws.Range("H19").Value = 6.25
$ws.Range("I19").Value = 624
$ws.Range("J19").Value = 618
$ws.Range("K19").Value = 0.970873786407
$ws.Range("L19").Value = 40.224719101123
$ws.Range("M19").Value = -28.358208955223
$ws.Range("N19").Value = -62.545018007202
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -66.666666666666
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = -20
$ws.Range("L20").Value = 21.212121212121
$ws.Range("M20").Value = 60
$ws.Range("N20").Value = -94.513031550068
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = -8.333333333333
$ws.Range("F21").Value = 138
$ws.Range("G21").Value = 142
$ws.Range("H21").Value = -2.81690140845
$ws.Range("I21").Value = 1049
$ws.Range("J21").Value = 1127
$ws.Range("K21").Value = -6.921029281277
$ws.Range("L21").Value = 22.690058479532
$ws.Range("M21").Value = -16.812053925456
$ws.Range("N21").Value = -75.306026365348
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -33.333333333333
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 40
$ws.Range("I22").Value = 57
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = -5
$ws.Range("L22").Value = 90
$ws.Range("M22").Value = 16.326530612244
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("C24").Value = 43
$ws.Range("D24").Value = 52
$ws.Range("E24").Value = -17.307692307692
$ws.Range("F24").Value = 187
$ws.Range("G24").Value = 213
$ws.Range("H24").Value = -12.206572769953
$ws.Range("I24").Value = 1335
$ws.Range("J24").Value = 1485
$ws.Range("K24").Value = -10.10101010101
$ws.Range("L24").Value = 24.766355140186
$ws.Range("M24").Value = 17.517605633802
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 53
$ws.Range("H25").Value = -15.094339622641
$ws.Range("I25").Value = 305
$ws.Range("J25").Value = 321
$ws.Range("K25").Value = -4.984423676012
$ws.Range("L25").Value = 19.607843137254
$ws.Range("M25").Value = 21.513944223107
$ws.Range("F26").Value = 3
$ws.Range("L26").Value = -19.047619047619
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -12.5
$ws.Range("I27").Value = 73
$ws.Range("J27").Value = 63
$ws.Range("K27").Value = 15.873015873015
$ws.Range("L27").Value = 4.285714285714
$ws.Range("D30").Value = 1
$ws.Range("G30").Value = 4
$ws.Range("J30").Value = 16
$ws.Range("K30").Value = -75
$ws.Range("L30").Value = -69.230769230769
